# Actualización automática 2025-08-07 15:30:08
# Corrects the "PIEDRA SINTERIZADA" sales figure for client
# "RENOVA&DISEÑA S.A." (advisor LOZANO MOLINA TITO) from 2227.24 to 1670.43
# across the three report sheets, and recomputes the dependent totals on
# the "CUMPLIMIENTO MENSUAL" summary sheet.

$wb = $excel.ActiveWorkbook

$oldValue = 2227.24
$newValue = 1670.43

# --- Sheet 1: "VENTAS POR GRUPO" -----------------------------------------
# Row 19 = RENOVA&DISEÑA S.A., column L = PIEDRA SINTERIZADA
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("L19").Value2 = $newValue

# --- Sheet 2: "VENTA MENSUAL" ---------------------------------------------
# Row 19 = RENOVA&DISEÑA S.A., column F = agosto
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F19").Value2 = $newValue

# Row 29 = TOTAL row, column F = agosto total
$wsVentaMensual.Range("F29").Value2 = $newValue

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" --------------------------------------
# Row 15 = PIEDRA SINTERIZADA group: VENTA (D), POR CUMPLIR (E), CUMPLIMIENTO (F)
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$presupuesto15 = $wsCumplimiento.Range("C15").Value2
$wsCumplimiento.Range("D15").Value2 = $newValue
$wsCumplimiento.Range("E15").Value2 = $presupuesto15 - $newValue
$wsCumplimiento.Range("F15").Value2 = $newValue / $presupuesto15

# Row 19 = TOTAL row: VENTA (D), POR CUMPLIR (E), CUMPLIMIENTO (F)
$presupuesto19 = $wsCumplimiento.Range("C19").Value2
$wsCumplimiento.Range("D19").Value2 = $newValue
$wsCumplimiento.Range("E19").Value2 = $presupuesto19 - $newValue
$wsCumplimiento.Range("F19").Value2 = $newValue / $presupuesto19

# Widen column F slightly (22 -> 23 characters) to match the refreshed report
# layout. Excel's ColumnWidth (COM) is expressed in "Normal style" characters
# and gets pixel-quantized (6px per character + 5px padding) before being
# written back out as the OOXML <col width="..."/> value, so 22.17 is the
# COM-side width that round-trips to an OOXML width of exactly 23.
$wsCumplimiento.Columns.Item(6).ColumnWidth = 22.17
